$wb = $excel.ActiveWorkbook

# --- Enable iterative calculation (Application-level settings) ---
# Mirrors calcPr iterate="1" iterateDelta="1.0000000000000001E-5" in the
# saved workbook.xml (iterateCount left at Excel's default of 100).
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.00001

# --- Update the ETLE exponent value and the sheet's remembered selection ---
$wsAbout = $wb.Worksheets.Item("About")
$wsETLE = $wb.Worksheets.Item("ETLE")

$wsETLE.Range("B2").Value = -5
$wsETLE.Range("B3").Select()

# Restore "About" as the active sheet/tab (selection on ETLE is merely the
# cached cursor position for that sheet, not the active tab).
$wsAbout.Activate()
